$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 716.3333
$ws.Range("I15").Value = 716.3333
$ws.Range("K15").Value = 2148.9999
$ws.Range("M15").Value = -1979.9999
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480
$ws.Range("H127").Value = 2445.8333
$ws.Range("I127").Value = 2445.8333
$ws.Range("K127").Value = 7337.499899999999
$ws.Range("M127").Value = -2377.499899999999
$ws.Range("H132").Value = 6099.6
$ws.Range("I132").Value = 7373.5
$ws.Range("K132").Value = 22120.5
$ws.Range("M132").Value = -19590.5
$ws.Range("H137").Value = 4459.75
$ws.Range("I137").Value = 4286
$ws.Range("K137").Value = 12858
$ws.Range("M137").Value = -10308
$ws.Range("H138").Value = 2928.4285
$ws.Range("I138").Value = 2749.5
$ws.Range("K138").Value = 8248.5
$ws.Range("M138").Value = -3108.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 879.6667
$ws.Range("I2").Value = 1269.5
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1269.5
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -1156.5
$ws.Range("N2").Value = -326
$ws.Range("H55").Value = 25000
$ws.Range("J55").Value = 25000
$ws.Range("L55").Value = 25000
$ws.Range("N55").Value = -25630
$ws.Range("H63").Value = 2154.6
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2154.6
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H116").Value = 879.6667
$ws.Range("I116").Value = 1269.5
$ws.Range("J116").Value = 100
$ws.Range("K116").Value = 1269.5
$ws.Range("L116").Value = 100
$ws.Range("M116").Value = 1024.5
$ws.Range("N116").Value = -4688

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 879.6667
$ws.Range("I3").Value = 1269.5
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 1269.5
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = -1155.5
$ws.Range("N3").Value = -328
$ws.Range("H82").Value = 19327.46
$ws.Range("I82").Value = 14584.111
$ws.Range("K82").Value = 14584.111
$ws.Range("M82").Value = -14201.111
$ws.Range("H85").Value = 19327.46
$ws.Range("I85").Value = 14584.111
$ws.Range("K85").Value = 14584.111
$ws.Range("M85").Value = -13258.111
$ws.Range("H86").Value = 6983.3335
$ws.Range("I86").Value = 5500
$ws.Range("J86").Value = 9950
$ws.Range("K86").Value = 5500
$ws.Range("L86").Value = 9950
$ws.Range("M86").Value = -4377
$ws.Range("N86").Value = -12196
$ws.Range("H89").Value = 6983.3335
$ws.Range("I89").Value = 5500
$ws.Range("J89").Value = 9950
$ws.Range("K89").Value = 27500
$ws.Range("L89").Value = 49750
$ws.Range("M89").Value = -21884
$ws.Range("N89").Value = -60982
$ws.Range("H99").Value = 2610.6924
$ws.Range("I99").Value = 2630.818
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2630.818
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -1132.818
$ws.Range("N99").Value = -5496

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 7165.75
$ws.Range("I25").Value = 4550
$ws.Range("K25").Value = 4550
$ws.Range("M25").Value = -4376
$ws.Range("H50").Value = 20070.215
$ws.Range("I50").Value = 20126.625
$ws.Range("J50").Value = 19995
$ws.Range("K50").Value = 20126.625
$ws.Range("L50").Value = 19995
$ws.Range("M50").Value = -19501.625
$ws.Range("N50").Value = -21245
$ws.Range("H70").Value = 21250
$ws.Range("J70").Value = 21250
$ws.Range("L70").Value = 21250
$ws.Range("N70").Value = -21880
$ws.Range("H73").Value = 21250
$ws.Range("J73").Value = 21250
$ws.Range("L73").Value = 21250
$ws.Range("N73").Value = -23434

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 3011.25
$ws.Range("I18").Value = 3465
$ws.Range("J18").Value = 1650
$ws.Range("K18").Value = 10395
$ws.Range("L18").Value = 4950
$ws.Range("M18").Value = -10226
$ws.Range("N18").Value = -5288
$ws.Range("H46").Value = 1371.125
$ws.Range("I46").Value = 994.8333
$ws.Range("J46").Value = 2500
$ws.Range("K46").Value = 2984.4999
$ws.Range("L46").Value = 7500
$ws.Range("M46").Value = -2893.4999
$ws.Range("N46").Value = -7682
$ws.Range("H103").Value = 214.8
$ws.Range("I103").Value = 23
$ws.Range("K103").Value = 69
$ws.Range("M103").Value = 810
$ws.Range("H109").Value = 889.5714
$ws.Range("I109").Value = 704.5
$ws.Range("K109").Value = 2113.5
$ws.Range("M109").Value = -1073.5
$ws.Range("H118").Value = 3000
$ws.Range("I118").Value = 3000
$ws.Range("K118").Value = 9000
$ws.Range("M118").Value = -7757
$ws.Range("H123").Value = 4497.75
$ws.Range("I123").Value = 4997
$ws.Range("J123").Value = 3000
$ws.Range("K123").Value = 14991
$ws.Range("L123").Value = 9000
$ws.Range("M123").Value = -12541
$ws.Range("N123").Value = -13900
$ws.Range("H124").Value = 4000
$ws.Range("J124").Value = 4000
$ws.Range("L124").Value = 12000
$ws.Range("N124").Value = -21820
$ws.Range("H125").Value = 7000
$ws.Range("J125").Value = 7000
$ws.Range("L125").Value = 21000
$ws.Range("N125").Value = -30840
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -10060
$ws.Range("H130").Value = 1030
$ws.Range("I130").Value = 1030
$ws.Range("K130").Value = 3090
$ws.Range("M130").Value = 1930

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19333.334
$ws.Range("J57").Value = 19200
$ws.Range("L57").Value = 19200
$ws.Range("N57").Value = -20840
$ws.Range("H122").Value = 1999.5
$ws.Range("I122").Value = 1999.5
$ws.Range("K122").Value = 5998.5
$ws.Range("M122").Value = -3548.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4249.25
$ws.Range("I16").Value = 4265.6665
$ws.Range("J16").Value = 4200
$ws.Range("K16").Value = 4265.6665
$ws.Range("L16").Value = 4200
$ws.Range("M16").Value = -4095.6665
$ws.Range("N16").Value = -4540
$ws.Range("H93").Value = 1071.9
$ws.Range("I93").Value = 964.875
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 964.875
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 283.125
$ws.Range("N93").Value = -3996
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 3624
$ws.Range("I132").Value = 2884.9092
$ws.Range("J132").Value = 4785.4287
$ws.Range("K132").Value = 8654.7276
$ws.Range("L132").Value = 14356.2861
$ws.Range("M132").Value = -6124.7276
$ws.Range("N132").Value = -19416.2861
$ws.Range("H136").Value = 2999
$ws.Range("I136").Value = 2999
$ws.Range("K136").Value = 8997
$ws.Range("M136").Value = -6447

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 15626
$ws.Range("J74").Value = 15626
$ws.Range("L74").Value = 15626
$ws.Range("N74").Value = -17498
$ws.Range("H77").Value = 15626
$ws.Range("J77").Value = 15626
$ws.Range("L77").Value = 46878
$ws.Range("N77").Value = -56238
$ws.Range("H103").Value = 16968
$ws.Range("J103").Value = 16968
$ws.Range("L103").Value = 16968
$ws.Range("N103").Value = -19312
